$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 2 of the "נתמכים" table with a new record.
# Cells are set in the order that makes new shared strings get
# registered with the same index ordering as the target workbook.
$ws.Range("B2").Value = "שפרינצק"
$ws.Range("F2").Value = "012-3456789"
$ws.Range("G2").Value = "987-6543210"
$ws.Range("H2").Value = "רווחה"
$ws.Range("A2").Value = "הוסרה מהגמח"
$ws.Range("I2").Value = "27 יולי 2023"
$ws.Range("J2").Value = "לא עונה לטלפון"

$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 1

$ws.Range("J2").Select()
